$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "109×4=436" "348×3=1044"
Replace-Text "740×9=6660" "792×6=4752"
Replace-Text "881×6=5286" "133×9=1197"
Replace-Text "204×4=816" "332×7=2324"
Replace-Text "905×6=5430" "669×2=1338"
Replace-Text "519×7=3633" "248×3=744"
Replace-Text "135×3=405" "184×4=736"
Replace-Text "639×7=4473" "785×4=3140"
Replace-Text "382×6=2292" "346×2=692"
Replace-Text "871×8=6968" "649×3=1947"
Replace-Text "122×3=366" "205×8=1640"
Replace-Text "431×7=3017" "590×4=2360"
Replace-Text "586×7=4102" "792×9=7128"
Replace-Text "158×3=474" "887×4=3548"
Replace-Text "607×7=4249" "599×5=2995"
Replace-Text "152×8=1216" "409×3=1227"
Replace-Text "465×6=2790" "601×9=5409"
Replace-Text "570×9=5130" "307×4=1228"
Replace-Text "598×5=2990" "548×2=1096"
Replace-Text "582×7=4074" "388×3=1164"
Replace-Text "154×4=616" "539×8=4312"
Replace-Text "752×8=6016" "766×5=3830"
Replace-Text "981×7=6867" "451×5=2255"
Replace-Text "276×8=2208" "193×4=772"
Replace-Text "818×2=1636" "422×7=2954"

Write-Output "Done applying replacements"
